# Update Work Week and Social Spending
#
# Refreshes the "GDP per Capita" data series for Albania on the Data sheet:
#  - revises the existing Data value (col E) for every year from 1870 to 2010
#  - appends six new rows, one per year, for 2011-2016

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Revise the existing yearly Data values (column E) ---
# The "Data" column stores these figures as text, even though they look like
# numbers. A plain `.Value = "711"` would get auto-converted to a Number, so
# each target cell is switched to Text format for the write and then the
# temporary format is cleared again so the cell keeps the workbook default
# (General/style 0) appearance.
$updates = @{
    2 = "711"
    22 = "953"
    32 = "1092"
    42 = "1243"
    45 = "1293"
    61 = "1476"
    82 = "1596"
    83 = "1666"
    84 = "1667"
    85 = "1736"
    86 = "1785"
    87 = "1882"
    88 = "1902"
    89 = "2023"
    90 = "2114"
    91 = "2201"
    92 = "2313"
    93 = "2332"
    94 = "2409"
    95 = "2491"
    96 = "2576"
    97 = "2670"
    98 = "2770"
    99 = "2876"
    100 = "2979"
    101 = "3080"
    102 = "3194"
    103 = "3322"
    104 = "3451"
    105 = "3623"
    106 = "3637"
    107 = "3649"
    108 = "3665"
    109 = "3681"
    110 = "3696"
    111 = "3716"
    112 = "3741"
    113 = "3762"
    114 = "3783"
    115 = "3805"
    116 = "3826"
    117 = "3846"
    118 = "3870"
    119 = "3894"
    120 = "3920"
    121 = "3948"
    122 = "3983"
    123 = "2942.10059831345"
    124 = "2841.76516466455"
    125 = "3137.0239335031"
    126 = "3399.71018272432"
    127 = "3666.65056914009"
    128 = "3965.68531531917"
    129 = "3526.98564071555"
    130 = "3873.37315205246"
    131 = "4432.04263745532"
    132 = "4808.47961481352"
    133 = "5285.80909643499"
    134 = "5608.96234460107"
    135 = "6004.64523899215"
    136 = "6419.87039572498"
    137 = "6858.46704735724"
    138 = "7347.33012304625"
    139 = "7866.15774715018"
    140 = "8522.12984128861"
    141 = "8859.47159720472"
    142 = "9222.97275207776"
}
foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    $cell.ClearFormats()
}

# --- Append the new rows for 2011-2016 ---
$newYears = @(2011, 2012, 2013, 2014, 2015, 2016)
$newData = @("9484", "9592", "9660", "9808", "10032", "10342")
$startRow = 143
for ($i = 0; $i -lt $newYears.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Albania"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $newYears[$i]
    $dataCell = $ws.Cells.Item($r, 5)
    $dataCell.NumberFormat = "@"
    $dataCell.Value = $newData[$i]
    $dataCell.ClearFormats()
}
